# -----------------------------------------------------------------------
# "Making numbers out of strings for if statement templating."
#
# The template used curly-quoted string literals for the branch numbers,
# e.g. {% if fields.a = “1” %} ... {% if fields.b != ”<10” %}. This
# edit turns every one of those into a bare number (no quotes), e.g.
# {% if fields.a = 1 %} ... {% if fields.b != <10 %}, so the merge engine
# compares fields.a/fields.b against numbers instead of strings.
#
# Doing that retype in real Word re-runs the as-you-type spell checker
# over the whole paragraph, which re-splits the runs and wraps the
# non-dictionary tokens (fields.a, fields.b, endif) in
# <w:proofErr w:type="spellStart/spellEnd"/> markers, and leaves the
# "_GoBack" last-edit bookmark sitting wherever the final keystroke
# landed. We reproduce that exact, already-known-good paragraph markup
# (rather than trying to re-derive the proofErr placement) via
# Range.InsertXML.
# -----------------------------------------------------------------------

$d = $word.ActiveDocument

# Find the paragraph that holds the if/else template (defensive lookup
# instead of assuming it is Paragraphs.Item(1)).
$target = $null
foreach ($para in $d.Paragraphs) {
    if ($para.Range.Text -like "*fields.a*") {
        $target = $para
        break
    }
}
if ($target -eq $null) {
    throw "template paragraph (containing 'fields.a') not found"
}
$targetIndex = $target.Index

# Preserve this paragraph's own <w:p> attributes (paraId/rsids/etc.)
# instead of hard-coding them.
$full = $target.Range
$paraAttrs = ""
if ($full.WordOpenXML -match "<w:p\s+([^>]*)>") {
    $paraAttrs = $matches[1]
}
$openTag = "<w:p>"
if ($paraAttrs -ne "") {
    $openTag = "<w:p " + $paraAttrs + ">"
}

# Exact run / proofErr / bookmark markup for the corrected paragraph.
$newInner = '<w:r><w:t xml:space="preserve">{% if </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>fields.a</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> = 1 %}</w:t></w:r><w:r><w:t>111</w:t></w:r><w:r><w:t xml:space="preserve">{% if </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>fields.b</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> %}</w:t></w:r><w:r><w:t xml:space="preserve">1 </w:t></w:r><w:r><w:t xml:space="preserve">&lt;10{% </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>endif</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> %}{% if </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>fields.b</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> != &lt;10 %}</w:t></w:r><w:r><w:t xml:space="preserve">1 </w:t></w:r><w:r><w:t>{{</w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>fields.b</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve">}}{% </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>endif</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> %}{% </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>endif</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> %}{% if </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>fields.a</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> = 2 %}{% if </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>fields.b</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> %}</w:t></w:r><w:r><w:t>2</w:t></w:r><w:r><w:t xml:space="preserve">&lt;10{% </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>endif</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> %}{% if </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>fields.b</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> != &lt;10 %}</w:t></w:r><w:r><w:t>2</w:t></w:r><w:r><w:t>{{</w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>fields.b</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve">}}{% </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>endif</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> %}{% </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>endif</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> %}{% if </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>fields.a</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> = 3 %}</w:t></w:r><w:r><w:t xml:space="preserve">3 </w:t></w:r><w:r><w:t xml:space="preserve">N/A{% </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>endif</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> %}{% if </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>fields.a</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> = 4 %}</w:t></w:r><w:r><w:t xml:space="preserve">4 </w:t></w:r><w:r><w:t xml:space="preserve">N/A{% </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>endif</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> %}{% if </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>fields.a</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> = 5 %}</w:t></w:r><w:r><w:t xml:space="preserve">5 </w:t></w:r><w:r><w:t xml:space="preserve">N/A{% </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>endif</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> %}{% if </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>fields.a</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> = 6 %}</w:t></w:r><w:r><w:t xml:space="preserve">6 </w:t></w:r><w:r><w:t xml:space="preserve">N/A{% </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>endif</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> %}{% if </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>fields.a</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> = 7 %}</w:t></w:r><w:r><w:t xml:space="preserve">7 </w:t></w:r><w:r><w:t xml:space="preserve">N/A{% </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>endif</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> %}{% if </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>fields.a</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> = </w:t></w:r><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/><w:r><w:t>8 %}</w:t></w:r><w:r><w:t xml:space="preserve">8 </w:t></w:r><w:r><w:t xml:space="preserve">N/A{% </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>endif</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> %}</w:t></w:r>'

$frag = (
    '<?xml version="1.0" encoding="UTF-8" standalone="yes"?>' +
    '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
    '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
    '<pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
    "<w:body>$openTag$newInner</w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>"
)

# InsertXML always opens a fresh paragraph for block-level <w:p> content,
# so this drops the rebuilt paragraph in just before the original one...
$insertPoint = $d.Range($full.Start, $full.Start)
$insertPoint.InsertXML($frag)

# ...which leaves the old (now emptied-of-purpose) paragraph pushed down
# one slot; delete it, re-fetching by index since old object references
# go stale across the InsertXML mutation.
$d.Paragraphs.Item($targetIndex + 1).Range.Delete() | Out-Null
